# Updates the cryptos price table (Price / Volume(1h) columns, plus two
# rows whose coin identity swapped position) to reflect the refreshed
# GitHub Actions scrape.
#
# Note: several "Price" values (column D) are plain decimal-looking
# strings (e.g. "1.00", "33.12") that Excel would otherwise auto-convert
# to numbers on assignment. Prefixing them with a leading apostrophe
# forces Excel to keep them as literal text, matching the source data
# (which stores every price as text, including multi-dot values like
# "34.224.73" that can't be parsed as numbers anyway).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '34.224.73'
$ws.Range("E2").Value = '  -0.75%  '
$ws.Range("D3").Value = '1.807.55'
$ws.Range("E3").Value = '  +1.14%  '
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = '  -0.20%  '
$ws.Range("D5").Value = "'223.62"
$ws.Range("E5").Value = '  +0.51%  '
$ws.Range("D6").Value = "'0.554"
$ws.Range("E6").Value = '  +0.11%  '
$ws.Range("E7").Value = '  -0.15%  '
$ws.Range("D8").Value = "'33.12"
$ws.Range("E8").Value = '  +3.42%  '
$ws.Range("D9").Value = "'0.289"
$ws.Range("E9").Value = '  +2.93%  '
$ws.Range("D10").Value = "'0.0719"
$ws.Range("E10").Value = '  +4.89%  '
$ws.Range("D11").Value = "'0.0929"
$ws.Range("E11").Value = '  -0.69%  '
$ws.Range("D12").Value = '2.064.79'
$ws.Range("E12").Value = '  +1.04%  '
$ws.Range("D13").Value = "'11.10"
$ws.Range("E13").Value = '  +2.07%  '
$ws.Range("D14").Value = '1.801.29'
$ws.Range("E14").Value = '  +0.75%  '
$ws.Range("D15").Value = "'0.634"
$ws.Range("E15").Value = '  +0.67%  '
$ws.Range("D16").Value = '34.238.10'
$ws.Range("E16").Value = '  -0.79%  '
$ws.Range("D17").Value = "'4.25"
$ws.Range("E17").Value = '  -0.63%  '
$ws.Range("D18").Value = "'68.88"
$ws.Range("E18").Value = '  +0.36%  '
$ws.Range("D19").Value = "'248.38"
$ws.Range("E19").Value = '  -2.03%  '
$ws.Range("E20").Value = '  +1.08%  '
$ws.Range("D21").Value = "'11.09"
$ws.Range("E21").Value = '  +6.21%  '
$ws.Range("D22").Value = "'1.00"
$ws.Range("E22").Value = '  -0.09%  '
$ws.Range("D23").Value = "'4.13"
$ws.Range("E23").Value = '  -0.54%  '
$ws.Range("E24").Value = '  +0.03%  '
$ws.Range("D25").Value = "'159.82"
$ws.Range("E25").Value = '  -0.42%  '
$ws.Range("D26").Value = "'16.65"
$ws.Range("E26").Value = '  +1.70%  '
$ws.Range("D27").Value = "'7.18"
$ws.Range("E27").Value = '  +1.52%  '
$ws.Range("D29").Value = "'0.999"
$ws.Range("D30").Value = "'0.0530"
$ws.Range("E30").Value = '  +2.85%  '
$ws.Range("D31").Value = "'3.75"
$ws.Range("E31").Value = '  +0.10%  '
$ws.Range("E32").Value = '  +1.51%  '
$ws.Range("D33").Value = "'3.55"
$ws.Range("E33").Value = '  +0.03%  '
$ws.Range("E34").Value = '  -0.70%  '
$ws.Range("D35").Value = '1.420.81'
$ws.Range("E35").Value = '  -1.03%  '
$ws.Range("E36").Value = '  +2.73%  '
$ws.Range("E37").Value = '  +0.81%  '
$ws.Range("E38").Value = '  -0.54%  '
$ws.Range("D39").Value = "'0.949"
$ws.Range("E39").Value = '  +3.76%  '
$ws.Range("D40").Value = "'80.99"
$ws.Range("E40").Value = '  -4.28%  '
$ws.Range("B41").Value = 'MXToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D41").Value = "'2.73"
$ws.Range("E41").Value = '  -2.22%  '
$ws.Range("B42").Value = 'HuobiToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D42").Value = "'2.36"
$ws.Range("E42").Value = '  +0.98%  '
$ws.Range("D43").Value = "'2.17"
$ws.Range("E43").Value = '  +4.70%  '
$ws.Range("D44").Value = "'5.97"
$ws.Range("E44").Value = '  -0.42%  '
$ws.Range("D45").Value = "'108.64"
$ws.Range("E45").Value = '  +4.94%  '
$ws.Range("E46").Value = '  +1.01%  '
$ws.Range("B47").Value = 'WEMIXToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D47").Value = "'1.05"
$ws.Range("E47").Value = '  -1.01%  '
$ws.Range("B48").Value = 'RocketPoolETH'
$ws.Range("C48").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D48").Value = '1.964.96'
$ws.Range("E48").Value = '  +0.86%  '
$ws.Range("D49").Value = "'12.14"
$ws.Range("E49").Value = '  +1.35%  '
$ws.Range("E50").Value = '  -0.11%  '
$ws.Range("E51").Value = '  +2.71%  '
